$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Correct the product short-name text (missing hyphen after "248") on both sheets.
$newName = "248-MS-EI-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"
$ws1.Range("B1").Value = $newName
$ws2.Range("B1").Value = $newName

# Move the selection on the input sheet off A6:B6 onto B1.
[void]$ws1.Range("B1").Select()

# Make the output sheet the active / selected tab, with B1 selected.
[void]$ws2.Activate()
[void]$ws2.Range("B1").Select()
